# The worksheet had a standalone section-header row ("grandes regiões e
# unidades da federação") at row 6 with no data of its own; the real
# "norte" row (with its data) started only at row 7. This was a data
# alignment bug: delete that spurious header row so every region label
# lines up with its correct data, shifting all subsequent rows (and their
# data) up by one. The former last row (previously row 38, "distrito
# federal") now lands on row 37, and row 38 no longer exists.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(6).Delete()
